# Update "想去人数" (want-to-go count) figures in the F column on the
# "展览" sheet and the mirrored "全部类型" sheet, matching the refreshed
# scrape output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 1095
$wsExpo.Range("F3").Value  = 4200
$wsExpo.Range("F5").Value  = 327
$wsExpo.Range("F8").Value  = 36
$wsExpo.Range("F10").Value = 124
$wsExpo.Range("F11").Value = 305
$wsExpo.Range("F12").Value = 234
$wsExpo.Range("F13").Value = 2910
$wsExpo.Range("F14").Value = 139
$wsExpo.Range("F15").Value = 1443

# --- Sheet "全部类型" (all categories, mirrors the expo rows shifted by one) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 1095
$wsAll.Range("F3").Value  = 4200
$wsAll.Range("F5").Value  = 327
$wsAll.Range("F9").Value  = 36
$wsAll.Range("F11").Value = 124
$wsAll.Range("F12").Value = 305
$wsAll.Range("F13").Value = 234
$wsAll.Range("F14").Value = 2910
$wsAll.Range("F15").Value = 139
$wsAll.Range("F16").Value = 1443
